# Auto commit at 2025-10-27  7:59:19.89
# Appends two new daily rows (2025-10-26) for the two charging stations
# ("四方坪站充电量(kw)" and "高岭站充电量(kw)") to the bottom of the data
# table on Sheet1, mirroring the existing per-hour layout (columns C:Z
# hold the 24 hourly charge readings, column A the date serial, column B
# the station name).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 112
$ws.Cells.Item(112, 1).Value = 45956
$ws.Cells.Item(112, 2).Value = "四方坪站充电量(kw)"
$ws.Cells.Item(112, 3).Value = 629.23
$ws.Cells.Item(112, 4).Value = 1388.5029999999997
$ws.Cells.Item(112, 5).Value = 277.02100000000002
$ws.Cells.Item(112, 6).Value = 462.41500000000002
$ws.Cells.Item(112, 7).Value = 640.7639999999999
$ws.Cells.Item(112, 8).Value = 746.50300000000004
$ws.Cells.Item(112, 9).Value = 447.74399999999997
$ws.Cells.Item(112, 10).Value = 276.74900000000002
$ws.Cells.Item(112, 11).Value = 163.73999999999998
$ws.Cells.Item(112, 12).Value = 154.79999999999998
$ws.Cells.Item(112, 13).Value = 256.77
$ws.Cells.Item(112, 14).Value = 323.23700000000008
$ws.Cells.Item(112, 15).Value = 749.7299999999999
$ws.Cells.Item(112, 16).Value = 1156.8309999999999
$ws.Cells.Item(112, 17).Value = 554.84500000000014
$ws.Cells.Item(112, 18).Value = 252.60599999999997
$ws.Cells.Item(112, 19).Value = 350.75700000000001
$ws.Cells.Item(112, 20).Value = 264.11400000000003
$ws.Cells.Item(112, 21).Value = 100.18
$ws.Cells.Item(112, 22).Value = 106.56
$ws.Cells.Item(112, 23).Value = 89.054000000000002
$ws.Cells.Item(112, 24).Value = 98.962000000000003
$ws.Cells.Item(112, 25).Value = 167.01999999999998
$ws.Cells.Item(112, 26).Value = 37.070999999999998

# Row 113
$ws.Cells.Item(113, 1).Value = 45956
$ws.Cells.Item(113, 2).Value = "高岭站充电量(kw)"
$ws.Cells.Item(113, 3).Value = 405.16400000000004
$ws.Cells.Item(113, 4).Value = 391.51800000000003
$ws.Cells.Item(113, 5).Value = 224.29399999999998
$ws.Cells.Item(113, 6).Value = 156.797
$ws.Cells.Item(113, 7).Value = 15.297000000000001
$ws.Cells.Item(113, 8).Value = 155.69600000000003
$ws.Cells.Item(113, 9).Value = 105.41800000000001
$ws.Cells.Item(113, 10).Value = 118.53999999999999
$ws.Cells.Item(113, 11).Value = 281.46100000000001
$ws.Cells.Item(113, 12).Value = 338.99699999999996
$ws.Cells.Item(113, 13).Value = 28.116
$ws.Cells.Item(113, 14).Value = 243.30300000000003
$ws.Cells.Item(113, 15).Value = 492.25299999999993
$ws.Cells.Item(113, 16).Value = 162.977
$ws.Cells.Item(113, 17).Value = 352.91200000000003
$ws.Cells.Item(113, 18).Value = 396.93200000000002
$ws.Cells.Item(113, 19).Value = 155.72099999999998
$ws.Cells.Item(113, 20).Value = 28.703000000000003
$ws.Cells.Item(113, 21).Value = 60.279000000000003
$ws.Cells.Item(113, 22).Value = 58.762999999999998
$ws.Cells.Item(113, 23).Value = 30.273
$ws.Cells.Item(113, 24).Value = 16.213999999999999
$ws.Cells.Item(113, 25).Value = 0
$ws.Cells.Item(113, 26).Value = 0


# Update the sheet's scroll position / active selection to match the
# author's final view of the worksheet after the new rows were added.
$ws.Range("G116").Select() | Out-Null
